$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 3079.5454
$ws.Range("J98").Value = 4550
$ws.Range("L98").Value = 4550
$ws.Range("N98").Value = -7546
# Row 122
$ws.Range("H122").Value = 3079.5454
$ws.Range("J122").Value = 4550
$ws.Range("L122").Value = 13650
$ws.Range("N122").Value = -18550
# Row 129
$ws.Range("H129").Value = 1323880.9
$ws.Range("J129").Value = 1684824.9
$ws.Range("L129").Value = 5054474.699999999
$ws.Range("N129").Value = -5064474.699999999
# Row 132
$ws.Range("H132").Value = 3243.077
$ws.Range("I132").Value = 3441.818
$ws.Range("J132").Value = 2150
$ws.Range("K132").Value = 10325.454
$ws.Range("L132").Value = 6450
$ws.Range("M132").Value = -7795.454000000002
$ws.Range("N132").Value = -11510
# Row 133
$ws.Range("H133").Value = 40143.84
$ws.Range("J133").Value = 40143.84
$ws.Range("L133").Value = 40143.84
$ws.Range("N133").Value = -50263.84
# Row 136
$ws.Range("H136").Value = 56722.5
$ws.Range("J136").Value = 56722.5
$ws.Range("L136").Value = 56722.5
$ws.Range("N136").Value = -66922.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 13239.288
$ws.Range("I32").Value = 6618.326
$ws.Range("J32").Value = 64000
$ws.Range("K32").Value = 6618.326
$ws.Range("L32").Value = 64000
$ws.Range("M32").Value = -6331.326
$ws.Range("N32").Value = -64574
# Row 61
$ws.Range("H61").Value = 1216.4375
$ws.Range("I61").Value = 897.7273
$ws.Range("K61").Value = 897.7273
$ws.Range("M61").Value = -685.7273
# Row 74
$ws.Range("H74").Value = 1077.9048
$ws.Range("I74").Value = 1105.9445
$ws.Range("K74").Value = 1105.9445
$ws.Range("M74").Value = -231.9445000000001
# Row 77
$ws.Range("H77").Value = 1077.9048
$ws.Range("I77").Value = 1105.9445
$ws.Range("K77").Value = 5529.7225
$ws.Range("M77").Value = -1161.7225
# Row 131
$ws.Range("H131").Value = 47000
$ws.Range("J131").Value = 47000
$ws.Range("L131").Value = 47000
$ws.Range("N131").Value = -57080
# Row 132
$ws.Range("H132").Value = 1792.2424
$ws.Range("I132").Value = 1339.1904
$ws.Range("K132").Value = 4017.5712
$ws.Range("M132").Value = -1487.5712
# Row 136
$ws.Range("H136").Value = 1216.4375
$ws.Range("I136").Value = 897.7273
$ws.Range("K136").Value = 2693.1819
$ws.Range("M136").Value = -143.1819

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 16894.402
$ws.Range("I134").Value = 1548.6744
$ws.Range("J134").Value = 44388.832
$ws.Range("K134").Value = 4646.023200000001
$ws.Range("L134").Value = 133166.496
$ws.Range("M134").Value = -2111.023200000001
$ws.Range("N134").Value = -138236.496

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3629.1785
$ws.Range("I31").Value = 4496.3335
$ws.Range("J31").Value = 2068.3
$ws.Range("K31").Value = 4496.3335
$ws.Range("L31").Value = 2068.3
$ws.Range("M31").Value = -4201.3335
$ws.Range("N31").Value = -2658.3
# Row 34
$ws.Range("H34").Value = 3629.1785
$ws.Range("I34").Value = 4496.3335
$ws.Range("J34").Value = 2068.3
$ws.Range("K34").Value = 4496.3335
$ws.Range("L34").Value = 2068.3
$ws.Range("M34").Value = -4294.3335
$ws.Range("N34").Value = -2472.3
# Row 122
$ws.Range("H122").Value = 1919.7142
$ws.Range("I122").Value = 1481
$ws.Range("J122").Value = 2504.6667
$ws.Range("K122").Value = 4443
$ws.Range("L122").Value = 7514.000100000001
$ws.Range("M122").Value = -1993
$ws.Range("N122").Value = -12414.0001
# Row 132
$ws.Range("H132").Value = 3489.9048
$ws.Range("I132").Value = 2299.111
$ws.Range("J132").Value = 4383
$ws.Range("K132").Value = 6897.333
$ws.Range("L132").Value = 13149
$ws.Range("M132").Value = -4367.333
$ws.Range("N132").Value = -18209
# Row 134
$ws.Range("H134").Value = 20001848
$ws.Range("I134").Value = 1884.7
$ws.Range("J134").Value = 100001700
$ws.Range("K134").Value = 5654.1
$ws.Range("L134").Value = 300005100
$ws.Range("M134").Value = -3119.1
$ws.Range("N134").Value = -300010170

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 32.81818
$ws.Range("I12").Value = 50
$ws.Range("J12").Value = 23
$ws.Range("K12").Value = 150
$ws.Range("L12").Value = 69
$ws.Range("M12").Value = 23
$ws.Range("N12").Value = -415
# Row 113
$ws.Range("H113").Value = 39109.31
$ws.Range("I113").Value = 966.6667
$ws.Range("J113").Value = 50552.1
$ws.Range("K113").Value = 2900.0001
$ws.Range("L113").Value = 151656.3
$ws.Range("M113").Value = -730.0001000000002
$ws.Range("N113").Value = -155996.3
# Row 131
$ws.Range("H131").Value = 5506426.5
$ws.Range("J131").Value = 915.4358999999999
$ws.Range("L131").Value = 2746.3077
$ws.Range("N131").Value = -12826.3077
# Row 138
$ws.Range("H138").Value = 2495.7144
$ws.Range("I138").Value = 622.5
$ws.Range("K138").Value = 1867.5
$ws.Range("M138").Value = 3272.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 9035.643
$ws.Range("I61").Value = 11949.9
$ws.Range("J61").Value = 1750
$ws.Range("K61").Value = 11949.9
$ws.Range("L61").Value = 1750
$ws.Range("M61").Value = -11747.9
$ws.Range("N61").Value = -2154
# Row 113
$ws.Range("H113").Value = 9035.643
$ws.Range("I113").Value = 11949.9
$ws.Range("J113").Value = 1750
$ws.Range("K113").Value = 11949.9
$ws.Range("L113").Value = 1750
$ws.Range("M113").Value = -9779.9
$ws.Range("N113").Value = -6090

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 3850
$ws.Range("J62").Value = 3430
$ws.Range("L62").Value = 3430
$ws.Range("N62").Value = -4678
# Row 65
$ws.Range("H65").Value = 3850
$ws.Range("J65").Value = 3430
$ws.Range("L65").Value = 17150
$ws.Range("N65").Value = -23390
# Row 132
$ws.Range("H132").Value = 1732.2174
$ws.Range("I132").Value = 1603.7693
$ws.Range("J132").Value = 1899.2
$ws.Range("K132").Value = 4811.3079
$ws.Range("L132").Value = 5697.6
$ws.Range("M132").Value = -2281.3079
$ws.Range("N132").Value = -10757.6
# Row 135
$ws.Range("H135").Value = 73292.14
$ws.Range("J135").Value = 73292.14
$ws.Range("L135").Value = 73292.14
$ws.Range("N135").Value = -83432.14
# Row 140
$ws.Range("H140").Value = 51856.375
$ws.Range("J140").Value = 51856.375
$ws.Range("L140").Value = 51856.375
$ws.Range("N140").Value = -62216.375
# Row 141
$ws.Range("H141").Value = 95800
$ws.Range("J141").Value = 95800
$ws.Range("L141").Value = 95800
$ws.Range("N141").Value = -106160
